# "Agregados comentarios y se quitaron los datos harcodeados"
#
# Adds the (previously hard-coded-elsewhere) card-holder data as real
# cells on the "users and pass" sheet, columns C:H, and makes that sheet
# the active/selected one (it was "Articles list" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users and pass")

# --- New data columns (C:E plain text, F:H forced to text format) ---
# Row 1
$ws.Range("C1").Value = "Luis Fernando Diaz"
$ws.Range("D1").Value = "Mexico"
$ws.Range("E1").Value = "Aguascalientes"

# Row 2
$ws.Range("C2").Value = "Jose Juan Vazquez"
$ws.Range("D2").Value = "Mexico"
$ws.Range("E2").Value = "Aguascalientes"

# Card number / expiry month / expiry year - stored as text
$ws.Range("F1:H2").NumberFormat = "@"

$ws.Range("F1").Value = "4573123447582740"
$ws.Range("F2").Value = "4573123447583445"

$ws.Range("G1").Value = "12"
$ws.Range("G2").Value = "11"

$ws.Range("H1").Value = "2024"
$ws.Range("H2").Value = "2025"

# Column F needs to be wide enough to show the full card number
$ws.Columns.Item(6).EntireColumn.AutoFit()

# Printable area / orientation touched on this sheet
$ws.PageSetup.Orientation = 1

# This sheet becomes the active tab (previously "Articles list" was)
$ws.Activate() | Out-Null
$ws.Range("H7").Select() | Out-Null
